$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 62
$ws.Range("H62").Value = 7562.727
$ws.Range("I62").Value = 6797.75
$ws.Range("J62").Value = 7999.857
$ws.Range("K62").Value = 6797.75
$ws.Range("L62").Value = 7999.857
$ws.Range("M62").Value = -6173.75
$ws.Range("N62").Value = -9247.857
# row 65
$ws.Range("H65").Value = 7562.727
$ws.Range("I65").Value = 6797.75
$ws.Range("J65").Value = 7999.857
$ws.Range("K65").Value = 33988.75
$ws.Range("L65").Value = 39999.285
$ws.Range("M65").Value = -30868.75
$ws.Range("N65").Value = -46239.285

$ws = $wb.Worksheets.Item("ARM")
# row 4
$ws.Range("H4").Value = 550
$ws.Range("I4").Value = 550
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 550
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -434
# row 5
$ws.Range("H5").Value = 175
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 150
$ws.Range("M5").Value = -88
$ws.Range("N5").Value = -374
# row 61
$ws.Range("H61").Value = 2256
$ws.Range("I61").Value = 2256
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2256
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2044
# row 102
$ws.Range("H102").Value = 1384.6364
$ws.Range("I102").Value = 970
$ws.Range("J102").Value = 1540.125
$ws.Range("K102").Value = 970
$ws.Range("L102").Value = 1540.125
$ws.Range("M102").Value = 652
$ws.Range("N102").Value = -4784.125
# row 124
$ws.Range("H124").Value = 62714.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 62714.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 62714.5
$ws.Range("N124").Value = -72534.5
# row 136
$ws.Range("H136").Value = 2256
$ws.Range("I136").Value = 2256
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6768
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4218

$ws = $wb.Worksheets.Item("BSM")
# row 4
$ws.Range("H4").Value = 175
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = -85
$ws.Range("N4").Value = -380
# row 94
$ws.Range("H94").Value = 1935.3572
$ws.Range("I94").Value = 2026.5385
$ws.Range("J94").Value = 750
$ws.Range("K94").Value = 2026.5385
$ws.Range("L94").Value = 750
$ws.Range("M94").Value = -1575.5385
$ws.Range("N94").Value = -1652
# row 105
$ws.Range("H105").Value = 4034.2856
$ws.Range("I105").Value = 2832.6843
$ws.Range("J105").Value = 15449.5
$ws.Range("K105").Value = 2832.6843
$ws.Range("L105").Value = 15449.5
$ws.Range("M105").Value = -1085.6843
$ws.Range("N105").Value = -18943.5

$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 79.2
$ws.Range("I7").Value = 60.142857
$ws.Range("J7").Value = 123.666664
$ws.Range("K7").Value = 60.142857
$ws.Range("L7").Value = 123.666664
$ws.Range("M7").Value = 52.857143
$ws.Range("N7").Value = -349.666664
# row 10
$ws.Range("H10").Value = 442.4
$ws.Range("I10").Value = 442.4
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 442.4
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -303.4
# row 62
$ws.Range("H62").Value = 69165.664
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 102248.5
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 102248.5
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -103496.5
# row 65
$ws.Range("H65").Value = 69165.664
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 102248.5
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 511242.5
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -517482.5
# row 99
$ws.Range("H99").Value = 14294.518
$ws.Range("I99").Value = 10698.7
$ws.Range("J99").Value = 16187.053
$ws.Range("K99").Value = 10698.7
$ws.Range("L99").Value = 16187.053
$ws.Range("M99").Value = -9200.700000000001
$ws.Range("N99").Value = -19183.053
# row 126
$ws.Range("H126").Value = 14294.518
$ws.Range("I126").Value = 10698.7
$ws.Range("J126").Value = 16187.053
$ws.Range("K126").Value = 32096.1
$ws.Range("L126").Value = 48561.159
$ws.Range("M126").Value = -29626.1
$ws.Range("N126").Value = -53501.159
# row 132
$ws.Range("H132").Value = 3630
$ws.Range("I132").Value = 1949.75
$ws.Range("J132").Value = 6990.5
$ws.Range("K132").Value = 5849.25
$ws.Range("L132").Value = 20971.5
$ws.Range("M132").Value = -3319.25
$ws.Range("N132").Value = -26031.5
# row 134
$ws.Range("H134").Value = 3199.8235
$ws.Range("I134").Value = 2606.3845
$ws.Range("J134").Value = 5128.5
$ws.Range("K134").Value = 7819.1535
$ws.Range("L134").Value = 15385.5
$ws.Range("M134").Value = -5284.1535
$ws.Range("N134").Value = -20455.5

$ws = $wb.Worksheets.Item("CUL")
# row 22
$ws.Range("H22").Value = 1462.5
$ws.Range("I22").Value = 4800
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 14400
$ws.Range("L22").Value = 1050
$ws.Range("M22").Value = -14231
$ws.Range("N22").Value = -1388
# row 27
$ws.Range("H27").Value = 1462.5
$ws.Range("I27").Value = 4800
$ws.Range("J27").Value = 350
$ws.Range("K27").Value = 14400
$ws.Range("L27").Value = 1050
$ws.Range("M27").Value = -14298
$ws.Range("N27").Value = -1254
# row 107
$ws.Range("H107").Value = 586.9429
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 610.3939
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 1831.1817
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -5671.1817

$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 849.94116
$ws.Range("I97").Value = 824.2
$ws.Range("J97").Value = 886.7143
$ws.Range("K97").Value = 824.2
$ws.Range("L97").Value = 886.7143
$ws.Range("M97").Value = -328.2
$ws.Range("N97").Value = -1878.7143
# row 107
$ws.Range("H107").Value = 1335
$ws.Range("I107").Value = 2164.1667
$ws.Range("J107").Value = 952.3077
$ws.Range("K107").Value = 2164.1667
$ws.Range("L107").Value = 952.3077
$ws.Range("M107").Value = -244.1667000000002
$ws.Range("N107").Value = -4792.3077
# row 132
$ws.Range("H132").Value = 4403.1665
$ws.Range("I132").Value = 1198
$ws.Range("J132").Value = 6005.75
$ws.Range("K132").Value = 3594
$ws.Range("L132").Value = 18017.25
$ws.Range("M132").Value = -1064
$ws.Range("N132").Value = -23077.25

$ws = $wb.Worksheets.Item("LTW")
# row 55
$ws.Range("H55").Value = 885
$ws.Range("I55").Value = 856.5714
$ws.Range("J55").Value = 984.5
$ws.Range("K55").Value = 856.5714
$ws.Range("L55").Value = 984.5
$ws.Range("M55").Value = -683.5714
$ws.Range("N55").Value = -1330.5
# row 122
$ws.Range("H122").Value = 7639.25
$ws.Range("I122").Value = 6685.6665
$ws.Range("J122").Value = 10500
$ws.Range("K122").Value = 20056.9995
$ws.Range("L122").Value = 31500
$ws.Range("M122").Value = -17606.9995
$ws.Range("N122").Value = -36400

$ws = $wb.Worksheets.Item("WVR")
# row 51
$ws.Range("H51").Value = 39025.332
$ws.Range("I51").Value = 23499.5
$ws.Range("J51").Value = 70077
$ws.Range("K51").Value = 23499.5
$ws.Range("L51").Value = 70077
$ws.Range("M51").Value = -22989.5
$ws.Range("N51").Value = -71097
# row 81
$ws.Range("H81").Value = 13113.429
$ws.Range("I81").Value = 11333.333
$ws.Range("J81").Value = 14448.5
$ws.Range("K81").Value = 22666.666
$ws.Range("L81").Value = 28897
$ws.Range("M81").Value = -21605.666
$ws.Range("N81").Value = -31019
# row 84
$ws.Range("H84").Value = 13113.429
$ws.Range("I84").Value = 11333.333
$ws.Range("J84").Value = 14448.5
$ws.Range("K84").Value = 113333.33
$ws.Range("L84").Value = 144485
$ws.Range("M84").Value = -108029.33
$ws.Range("N84").Value = -155093
# row 92
$ws.Range("H92").Value = 122500
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 122500
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 122500
$ws.Range("N92").Value = -127492
